# Weekly data update: a new day's record is inserted at row 517, pushing the
# existing rows 517-585 down by one position (row 586 is newly created to hold
# what used to be the last row, 585). Only the data columns that vary row to
# row (D, J, K, L, M, O, P) need to be shifted - the remaining columns
# (A, B, C, E, F, G, H, I, N, Q, R) are constant for every row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 517
$lastRow  = 585
$newLastRow = 586

# 1) Snapshot the current (pre-edit) values of the columns that move, for
#    every row that will be shifted down (517..585). D is read with Value2
#    so we get the raw date serial number rather than a DateTime object -
#    this avoids Excel silently stamping a brand-new (default) date style
#    onto row 586 further down, which would leave an unused number format
#    behind in styles.xml.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rec = @{}
    $rec["D"] = $ws.Cells.Item($r, 4).Value2()
    $rec["J"] = $ws.Cells.Item($r, 10).Value()
    $rec["K"] = $ws.Cells.Item($r, 11).Value()
    $rec["L"] = $ws.Cells.Item($r, 12).Value()
    $rec["M"] = $ws.Cells.Item($r, 13).Value()
    $rec["O"] = $ws.Cells.Item($r, 15).Value()
    $rec["P"] = $ws.Cells.Item($r, 16).Value()
    $snapshot[$r] = $rec
}

# 2) Shift rows 518..586 down: new row N gets the pre-edit content of row N-1.
for ($r = $lastRow + 1; $r -ge $firstRow + 1; $r--) {
    $src = $snapshot[$r - 1]
    $ws.Cells.Item($r, 4).Value  = $src["D"]
    $ws.Cells.Item($r, 10).Value = $src["J"]
    $ws.Cells.Item($r, 11).Value = $src["K"]
    $ws.Cells.Item($r, 12).Value = $src["L"]
    $ws.Cells.Item($r, 13).Value = $src["M"]
    $ws.Cells.Item($r, 15).Value = $src["O"]
    $ws.Cells.Item($r, 16).Value = $src["P"]
}

# 3) Row 517 becomes the new record (its D and J values change; K, L, M, O, P
#    keep the values that were already there).
$ws.Cells.Item($firstRow, 4).Value  = 45154
$ws.Cells.Item($firstRow, 10).Value = 500

# 4) Row 586 is a brand new row - fill in the columns that are constant
#    across the whole sheet (these were never touched above because row 586
#    did not exist before).
$ws.Cells.Item($newLastRow, 1).Value  = 8
$ws.Cells.Item($newLastRow, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newLastRow, 3).Value  = "Coquimbo"
$ws.Cells.Item($newLastRow, 5).Value  = 4
$ws.Cells.Item($newLastRow, 6).Value  = 100114013
$ws.Cells.Item($newLastRow, 7).Value  = "Zanahoria"
$ws.Cells.Item($newLastRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newLastRow, 9).Value  = "Primera"
$ws.Cells.Item($newLastRow, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item($newLastRow, 17).Value = 20
$ws.Cells.Item($newLastRow, 18).Value = "Hortaliza"

# Make sure the new date cell (D586) renders/stores using the same date
# format as the rest of column D.
$ws.Cells.Item($newLastRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
